$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at Y (shifting the existing "c3"/"nota_iniciativa" columns
# to the right) and populate it with the new "l1" field.
$ws.Columns("Y").Insert()

$ws.Range("Y1").Value = "l1"

for ($r = 2; $r -le 62; $r++) {
    $ws.Cells.Item($r, 25).Value = 0
}
